$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1
$ws.Range("B1").Value = "Rapport de masculinité"

# Update data value in B2
$ws.Range("B2").Value = "122,7"

# Remove column C entirely (was Population de l'année 2019_1 / 5761)
$ws.Range("C1:C2").EntireColumn.Delete()
